$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics after the new trade closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1399.77   # Current Capital
$summary.Range("B4").Value = -0.23     # Total P&L $
$summary.Range("B5").Value = -1.53     # Total P&L %
$summary.Range("B6").Value = 3         # Total Trades
$summary.Range("B8").Value = 3         # Losing Trades

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 99.77      # Capital
$status.Range("D5").Value = 3          # Trades
$status.Range("E5").Value = -0.23      # P&L $
$status.Range("F5").Value = -0.23      # P&L %

# ---------------------------------------------------------------------------
# Append the newly-closed Trade #3 to both the "All Trades" and
# "MarketMaking" sheets (row 4)
# ---------------------------------------------------------------------------
$tradeSheets = @("All Trades", "MarketMaking")

foreach ($sheetName in $tradeSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A4").Value = 3

    # Dates such as "2026-02-17" are auto-detected by Excel and silently
    # converted to a date serial number. Force the cell to Text first so
    # the literal string is preserved, then drop back to the Normal style
    # so no stray number-format is left behind on the cell.
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2026-02-17"
    $ws.Range("B4").Style = "Normal"

    $ws.Range("C4").Value = "19:55:44"
    $ws.Range("D4").Value = "MarketMaking"
    $ws.Range("E4").Value = "DOWN"
    $ws.Range("F4").Value = 0.43
    $ws.Range("G4").Value = 0.42
    $ws.Range("H4").Value = "CLOSED"
    $ws.Range("I4").Value = -2.3256
    $ws.Range("J4").Value = -0.01
    $ws.Range("K4").Value = 99.77
    $ws.Range("L4").Value = 0
    $ws.Range("M4").Value = 0
    $ws.Range("N4").Value = 0.6
    $ws.Range("O4").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P4").Value = "early_exit"
    $ws.Range("Q4").Value = 0.13
}
